$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1950.3903
$ws.Range("I33").Value = 1710.1052
$ws.Range("K33").Value = 1710.1052
$ws.Range("M33").Value = -1481.1052
$ws.Range("H74").Value = 4490
$ws.Range("I74").Value = 4490
$ws.Range("K74").Value = 4490
$ws.Range("M74").Value = -3554
$ws.Range("H76").Value = 9349.645500000001
$ws.Range("I76").Value = 9312.583000000001
$ws.Range("J76").Value = 9476.714
$ws.Range("K76").Value = 9312.583000000001
$ws.Range("L76").Value = 9476.714
$ws.Range("M76").Value = -8997.583000000001
$ws.Range("N76").Value = -10106.714
$ws.Range("H77").Value = 4490
$ws.Range("I77").Value = 4490
$ws.Range("K77").Value = 22450
$ws.Range("M77").Value = -17770
$ws.Range("H79").Value = 9349.645500000001
$ws.Range("I79").Value = 9312.583000000001
$ws.Range("J79").Value = 9476.714
$ws.Range("K79").Value = 9312.583000000001
$ws.Range("L79").Value = 9476.714
$ws.Range("M79").Value = -8220.583000000001
$ws.Range("N79").Value = -11660.714
$ws.Range("H80").Value = 92537.27
$ws.Range("I80").Value = 1591.2
$ws.Range("K80").Value = 4773.6
$ws.Range("M80").Value = -3775.6
$ws.Range("H83").Value = 92537.27
$ws.Range("I83").Value = 1591.2
$ws.Range("K83").Value = 14320.8
$ws.Range("M83").Value = -9328.800000000001
$ws.Range("H86").Value = 5240.08
$ws.Range("I86").Value = 4964.7896
$ws.Range("J86").Value = 6111.8335
$ws.Range("K86").Value = 4964.7896
$ws.Range("L86").Value = 6111.8335
$ws.Range("M86").Value = -3841.7896
$ws.Range("N86").Value = -8357.833500000001
$ws.Range("H87").Value = 131563.5
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("H89").Value = 5240.08
$ws.Range("I89").Value = 4964.7896
$ws.Range("J89").Value = 6111.8335
$ws.Range("K89").Value = 24823.948
$ws.Range("L89").Value = 30559.1675
$ws.Range("M89").Value = -19207.948
$ws.Range("N89").Value = -41791.1675
$ws.Range("H90").Value = 131563.5
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("H92").Value = 836.4
$ws.Range("I92").Value = 659.7692
$ws.Range("K92").Value = 659.7692
$ws.Range("M92").Value = 588.2308
$ws.Range("H97").Value = 2583.5
$ws.Range("I97").Value = 394
$ws.Range("J97").Value = 3021.4
$ws.Range("K97").Value = 1182
$ws.Range("L97").Value = 9064.200000000001
$ws.Range("M97").Value = -686
$ws.Range("N97").Value = -10056.2
$ws.Range("H104").Value = 488
$ws.Range("I104").Value = 488
$ws.Range("K104").Value = 1464
$ws.Range("M104").Value = 283
$ws.Range("H116").Value = 8666.200000000001
$ws.Range("J116").Value = 8832.75
$ws.Range("L116").Value = 8832.75
$ws.Range("N116").Value = -15716.75
$ws.Range("H121").Value = 4441.3335
$ws.Range("J121").Value = 4441.3335
$ws.Range("L121").Value = 13324.0005
$ws.Range("N121").Value = -16818.0005
$ws.Range("H138").Value = 3334.634
$ws.Range("I138").Value = 1521.5
$ws.Range("J138").Value = 3774.182
$ws.Range("K138").Value = 4564.5
$ws.Range("L138").Value = 11322.546
$ws.Range("M138").Value = 575.5
$ws.Range("N138").Value = -21602.546
$ws.Range("M87").ClearContents()
$ws.Range("M90").ClearContents()

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2888.4614
$ws.Range("I74").Value = 2754.1667
$ws.Range("K74").Value = 2754.1667
$ws.Range("M74").Value = -1880.1667
$ws.Range("H77").Value = 2888.4614
$ws.Range("I77").Value = 2754.1667
$ws.Range("K77").Value = 13770.8335
$ws.Range("M77").Value = -9402.833500000001

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3628.875
$ws.Range("I105").Value = 4408.5
$ws.Range("J105").Value = 2849.25
$ws.Range("K105").Value = 4408.5
$ws.Range("L105").Value = 2849.25
$ws.Range("M105").Value = -2661.5
$ws.Range("N105").Value = -6343.25

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 80156
$ws.Range("J109").Value = 80156
$ws.Range("L109").Value = 80156
$ws.Range("N109").Value = -82236
$ws.Range("H134").Value = 6019.9844
$ws.Range("I134").Value = 5735.415
$ws.Range("J134").Value = 7391.091
$ws.Range("K134").Value = 17206.245
$ws.Range("L134").Value = 22173.273
$ws.Range("M134").Value = -14671.245
$ws.Range("N134").Value = -27243.273
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1931.7333
$ws.Range("I117").Value = 687.75
$ws.Range("K117").Value = 2063.25
$ws.Range("M117").Value = 1378.75

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("H122").Value = 3273.0833
$ws.Range("I122").Value = 3110.375
$ws.Range("J122").Value = 3598.5
$ws.Range("K122").Value = 9331.125
$ws.Range("L122").Value = 10795.5
$ws.Range("M122").Value = -6881.125
$ws.Range("N122").Value = -15695.5
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 18937.8
$ws.Range("I7").Value = 20426.916
$ws.Range("J7").Value = 12981.333
$ws.Range("K7").Value = 20426.916
$ws.Range("L7").Value = 12981.333
$ws.Range("M7").Value = -20314.916
$ws.Range("N7").Value = -13205.333
$ws.Range("H22").Value = 3418.6924
$ws.Range("I22").Value = 3327.8333
$ws.Range("J22").Value = 3496.5715
$ws.Range("K22").Value = 3327.8333
$ws.Range("L22").Value = 3496.5715
$ws.Range("M22").Value = -3032.8333
$ws.Range("N22").Value = -4086.5715
$ws.Range("H27").Value = 3418.6924
$ws.Range("I27").Value = 3327.8333
$ws.Range("J27").Value = 3496.5715
$ws.Range("K27").Value = 3327.8333
$ws.Range("L27").Value = 3496.5715
$ws.Range("M27").Value = -3220.8333
$ws.Range("N27").Value = -3710.5715
$ws.Range("H55").Value = 499.81818
$ws.Range("I55").Value = 543
$ws.Range("J55").Value = 424.25
$ws.Range("K55").Value = 543
$ws.Range("L55").Value = 424.25
$ws.Range("M55").Value = -370
$ws.Range("N55").Value = -770.25
$ws.Range("H61").Value = 2986.9
$ws.Range("I61").Value = 2431.8333
$ws.Range("J61").Value = 3819.5
$ws.Range("K61").Value = 2431.8333
$ws.Range("L61").Value = 3819.5
$ws.Range("M61").Value = -2229.8333
$ws.Range("N61").Value = -4223.5
$ws.Range("H97").Value = 52758
$ws.Range("J97").Value = 52758
$ws.Range("L97").Value = 52758
$ws.Range("N97").Value = -54740
$ws.Range("H113").Value = 2986.9
$ws.Range("I113").Value = 2431.8333
$ws.Range("J113").Value = 3819.5
$ws.Range("K113").Value = 2431.8333
$ws.Range("L113").Value = 3819.5
$ws.Range("M113").Value = -261.8332999999998
$ws.Range("N113").Value = -8159.5
$ws.Range("H126").Value = 18937.8
$ws.Range("I126").Value = 20426.916
$ws.Range("J126").Value = 12981.333
$ws.Range("K126").Value = 61280.74800000001
$ws.Range("L126").Value = 38943.999
$ws.Range("M126").Value = -58810.74800000001
$ws.Range("N126").Value = -43883.999
$ws.Range("H136").Value = 1827.3914
$ws.Range("I136").Value = 1743.6154
$ws.Range("J136").Value = 1936.3
$ws.Range("K136").Value = 5230.8462
$ws.Range("L136").Value = 5808.9
$ws.Range("M136").Value = -2680.8462
$ws.Range("N136").Value = -10908.9

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2683.5417
$ws.Range("I107").Value = 2491.1052
$ws.Range("J107").Value = 3414.8
$ws.Range("K107").Value = 7473.3156
$ws.Range("L107").Value = 10244.4
$ws.Range("M107").Value = -5553.3156
$ws.Range("N107").Value = -14084.4
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("H136").Value = 2695.111
$ws.Range("I136").Value = 2429.475
$ws.Range("K136").Value = 7288.424999999999
$ws.Range("M136").Value = -4738.424999999999
$ws.Range("N109").ClearContents()
